# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet (fund position detail) right after the
# "总计" summary sheet -- pushing 2022-Q3 .. 2021-Q2 down by one tab -- and
# adds the matching roll-up row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert the 2022-Q4 totals as the new first
#    data row, shifting the existing quarters (and their running index
#    in column A) down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Give the new row 8 / A8 the same look (border+bold+center style) as the
# rest of column A before we fill it in, by cloning the formatting of the
# row above.
$summary.Cells.Item(7, 1).Copy($summary.Cells.Item(8, 1))

$summaryRows = @(
    @("2022-Q4", 6,  0.18),
    @("2022-Q3", 13, 0.33),
    @("2022-Q2", 10, 2.03),
    @("2022-Q1", 2,  0.11),
    @("2021-Q4", 12, 4.79),
    @("2021-Q3", 3,  0.67),
    @("2021-Q2", 6,  0.25)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $i
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
}

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right before "2022-Q3" (i.e.
#    right after "总计"). Clone the layout/styling of the existing
#    "2022-Q3" sheet (header row + per-column formatting) so the new
#    sheet matches the look of its siblings, then overwrite the cell
#    values with the 2022-Q4 figures.
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$ws = $wb.Worksheets.Add($beforeSheet)
$ws.Name = "2022-Q4"

$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Range("B1:H7").Copy($ws.Range("B1:H7"))
$q3.Range("A2:A7").Copy($ws.Range("A2:A7"))

# Columns B, D..G hold figures that look numeric ("1.29", "0.0528", fund
# codes with leading zeros, ...) but must stay text, exactly like the
# other quarter sheets -- force Text format before writing them. Column C
# (fund name) never looks numeric, so it doesn't need this treatment.
$ws.Range("B2:B7").NumberFormat = "@"
$ws.Range("D2:G7").NumberFormat = "@"

$fundRows = @(
    @("016283", "华泰柏瑞积极优选股票C",       "1.29", "89.39", "4.09", "0.0528", 3),
    @("001097", "华泰柏瑞积极优选股票A",       "1.23", "89.39", "4.09", "0.0503", 3),
    @("014839", "兴银碳中和主题混合C",         "0.64", "92.17", "4.69", "0.0300", 4),
    @("014838", "兴银碳中和主题混合A",         "0.53", "92.17", "4.69", "0.0249", 4),
    @("009937", "东方欣益一年持有期偏债混合A", "1.93", "22.19", "0.83", "0.0160", 5),
    @("009938", "东方欣益一年持有期偏债混合C", "0.31", "22.19", "0.83", "0.0026", 5)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $ws.Cells.Item($r, 1).Value = $i
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
}
